# Insert a new data row at row 692 (pushing the existing rows 692:782 down
# to 693:783) and populate it with the new price-report record for
# "Macroferia Regional de Talca" (Uva / Red Globe / Segunda).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 692; this shifts rows 692-782
# down to 693-783 and grows the sheet dimension to A1:T783 automatically.
$ws.Rows.Item(692).Insert()

# Fill in the new record in the now-empty row 692.
$ws.Range("A692").Value = 5
$ws.Range("B692").Value = 'Macroferia Regional de Talca'
$ws.Range("C692").Value = 'Maule'
$ws.Range("D692").Value = 45127
$ws.Range("E692").Value = 7
$ws.Range("F692").Value = 'Fruta'
$ws.Range("G692").Value = 100109
$ws.Range("H692").Value = 'Uva'
$ws.Range("I692").Value = 100109001
$ws.Range("J692").Value = 'Uva'
$ws.Range("K692").Value = 'Red Globe'
$ws.Range("L692").Value = 'Segunda'
$ws.Range("M692").Value = 300
$ws.Range("N692").Value = 11000
$ws.Range("O692").Value = 11000
$ws.Range("P692").Value = 11000
$ws.Range("Q692").Value = '$/bandeja 8 kilos'
$ws.Range("R692").Value = "Región de O'Higgins"
$ws.Range("S692").Value = 1375
$ws.Range("T692").Value = 8
